$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = "RAVI KUMAR"
$ws.Range("B7").Value = "2345 6789 1234"
$ws.Range("C7").Value = "'12/05/1998"
$ws.Range("D7").Value = "Male"
